$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B18").Value = "est_propF"
$ws.Range("C18").Value = "Is sex ration F/(M+F) to be included in the likelihood (assumed normal); 0 = no, 1 = use annual average across ages (uses 2nd age in propF data), 2 = age, and year specific (TBD)"

$ws.Range("B19").Value = "propF_sigma"
$ws.Range("C19").Value = "Initial value or fixed value for sd of normal likelihood for sex ration. Not yet able to estimate."

$ws.Range("B20").Value = "fleet_control"
$ws.Range("C20").Value = "Survey and fishery data specifications"

$ws.Range("B21").Value = "srv_biom"
$ws.Range("C21").Value = "Survey index in weight (kg) or numbers data"

$ws.Range("B22").Value = "srv_emp_sel"
$ws.Range("C22").Value = "Empirical selectivity for surveys (leave empty if not used)"

$ws.Range("B23").Value = "comp"
$ws.Range("C23").Value = "Survey/fishery age or length composition data"

$ws.Range("B24").Value = "fsh_biom"
$ws.Range("C24").Value = "Total catch in weight (kg) or numbers data"

$ws.Range("B25").Value = "fsh_emp_sel"
$ws.Range("C25").Value = "Empirical selectivity for fisheries (leave empty if not used)"

$ws.Range("A1").Select()
$ws.Range("C7").Select()
